$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new source row (row 7) to the worksheet
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Μαθήματα με τον Μακσίμ Κισιλιέρ"

# Extend the Excel Table ("Tabla1") to include the newly added row
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:C7"))

# Update the active selection to match the post-edit state (A7)
$ws.Range("A7").Select()
